# Update the workbook metadata and the Elements table to reflect the
# LinuxForHealth re-branding / republish of the StructureDefinition.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# --- Metadata sheet updates -------------------------------------------------
# URL
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/industry-classification"
# Version
$wsMeta.Range("B3").Value = "8.0.0"
# Date
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
# Publisher
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates --------------------------------------------------
# The Extension.url row's Fixed Value mirrors the StructureDefinition's own
# canonical URL, so it is updated in lockstep with the Metadata URL above.
$wsElem.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/industry-classification"

# Clear the Constraint(s) value for the base "Extension" row (row 2); that
# constraint text now only appears on the "Extension.extension" row.
$wsElem.Range("AI2").Value = ""
